# PNAD 2009 - seguranca / cv121203a
# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had a stray row (row 6) whose label was the shared string
# "grandes regiões e unidades da federação" but which never got any of the
# B:I numeric data next to it - the real per-region figures all started one
# row too low (row 7 held the "norte" data misaligned under the blank
# "grandes regiões..." row, etc.), leaving a trailing empty-label row
# (old row 37, "goiás") with no matching data as well.
#
# Fix: remove that spurious row entirely. Excel shifts every row below it
# up by one, which re-aligns each region's label with its correct data and
# drops the now-unused trailing row, shrinking the sheet from A1:I37 to
# A1:I36 and removing the now-unreferenced shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").EntireRow.Delete()
